$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Pre-apply the existing cell style (Arial / General, style index 1, taken
#    from an already-styled cell) to every new cell we are about to populate,
#    BEFORE writing any values. This reuses the workbook's existing style
#    record instead of Excel fabricating brand-new (unused) font/xf entries.
# ---------------------------------------------------------------------------
$ws.Range("G4").Copy()
$ws.Range("A5:H14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1:H4").PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) New user rows 5-8 (columns A-G only) - first batch of new people.
# ---------------------------------------------------------------------------

# Row 5: Diego Ramírez
$ws.Range("A5").Value = "Diego"
$ws.Range("B5").Value = "Ramírez"
$ws.Range("C5").Value = "diego.ramirez91@example.com"
$ws.Range("D5").Value = 3012233445
$ws.Range("E5").Value = "Dieg0!Test"
$ws.Range("F5").Value = "Dieg0!Test"
$ws.Range("G5").Value = "No"

# Row 6: Mariana Torres
$ws.Range("A6").Value = "Mariana"
$ws.Range("B6").Value = "Torres"
$ws.Range("C6").Value = "mariana.torres88@example.com"
$ws.Range("D6").Value = 3109876543
$ws.Range("E6").Value = "Mari#2025"
$ws.Range("F6").Value = "Mari#2025"
$ws.Range("G6").Value = "Yes"

# Row 7: Felipe Rodríguez
$ws.Range("A7").Value = "Felipe"
$ws.Range("B7").Value = "Rodríguez"
$ws.Range("C7").Value = "f.rodriguez@example.com"
$ws.Range("D7").Value = 3123344556
$ws.Range("E7").Value = "FeliPass_88"
$ws.Range("F7").Value = "FeliPass_88"
$ws.Range("G7").Value = "Yes"

# Row 8: Camila Herrera
$ws.Range("A8").Value = "Camila"
$ws.Range("B8").Value = "Herrera"
$ws.Range("C8").Value = "camila.h@example.com"
$ws.Range("D8").Value = 3135566778
$ws.Range("E8").Value = "CamH#321"
$ws.Range("F8").Value = "CamH#321"
$ws.Range("G8").Value = "No"

# ---------------------------------------------------------------------------
# 3) New column H ("used") - header plus the first "true "/"false" markers.
#    "false" (and "true" without trailing space) look like booleans to
#    Excel's auto-detection, so they are entered as formulas producing a
#    text string and then converted to a static value, which keeps them as
#    genuine shared-string text instead of boolean cells.
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "used"
$ws.Range("H2").Value = "true "

$ws.Range("H6").Formula = "=""false"""
$ws.Range("H6").Copy()
$ws.Range("H6").PasteSpecial(-4163)       # xlPasteValues
$excel.CutCopyMode = 0

# Remaining "true " rows for the already-existing + Diego rows.
$ws.Range("H3").Value = "true "
$ws.Range("H4").Value = "true "
$ws.Range("H5").Value = "true "

# Remaining "false" rows (7 and 8); done via the same formula trick so the
# cells stay text (t="s") instead of becoming boolean cells.
$ws.Range("H7").Formula = "=""false"""
$ws.Range("H8").Formula = "=""false"""
$ws.Range("H7:H8").Copy()
$ws.Range("H7:H8").PasteSpecial(-4163)    # xlPasteValues
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) New user rows 9-14 (all columns A-H).
# ---------------------------------------------------------------------------

# Row 9: Andrés Martínez
$ws.Range("A9").Value = "Andrés"
$ws.Range("B9").Value = "Martínez"
$ws.Range("C9").Value = "andres.mtz@example.com"
$ws.Range("D9").Value = 3147788990
$ws.Range("E9").Value = "AndMart2025!"
$ws.Range("F9").Value = "AndMart2025!"
$ws.Range("G9").Value = "Yes"
$ws.Range("H9").Formula = "=""false"""

# Row 10: Valentina Salazar
$ws.Range("A10").Value = "Valentina"
$ws.Range("B10").Value = "Salazar"
$ws.Range("C10").Value = "v.salazar@example.com"
$ws.Range("D10").Value = 3056677889
$ws.Range("E10").Value = "ValeTest#1"
$ws.Range("F10").Value = "ValeTest#1"
$ws.Range("G10").Value = "No"
$ws.Range("H10").Formula = "=""false"""

# Row 11: Tomás Ortega (telephone stored as TEXT, like rows 12-14)
$ws.Range("A11").Value = "Tomás"
$ws.Range("B11").Value = "Ortega"
$ws.Range("C11").Value = "tomas.ortega@example.com"
$ws.Range("D11").Formula = "=""3004455667"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)       # xlPasteValues
$excel.CutCopyMode = 0
$ws.Range("E11").Value = "T0mPass!23"
$ws.Range("F11").Value = "T0mPass!23"
$ws.Range("G11").Value = "Yes"
$ws.Range("H11").Formula = "=""false"""

# Row 12: Juliana Núñez (telephone stored as TEXT)
$ws.Range("A12").Value = "Juliana"
$ws.Range("B12").Value = "Núñez"
$ws.Range("C12").Value = "juliana.nunez@example.com"
$ws.Range("D12").Formula = "=""3112233445"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("E12").Value = "Juli_Nz2025"
$ws.Range("F12").Value = "Juli_Nz2025"
$ws.Range("G12").Value = "Yes"
$ws.Range("H12").Formula = "=""false"""

# Row 13: Samuel Pérez (telephone stored as TEXT)
$ws.Range("A13").Value = "Samuel"
$ws.Range("B13").Value = "Pérez"
$ws.Range("C13").Value = "samuel.perez@example.com"
$ws.Range("D13").Formula = "=""3185566770"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("E13").Value = "SamP!2024"
$ws.Range("F13").Value = "SamP!2024"
$ws.Range("G13").Value = "No"
$ws.Range("H13").Formula = "=""false"""

# Row 14: Daniela Morales (telephone stored as TEXT)
$ws.Range("A14").Value = "Daniela"
$ws.Range("B14").Value = "Morales"
$ws.Range("C14").Value = "daniela.morales@example.com"
$ws.Range("D14").Formula = "=""3213344556"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("E14").Value = "D@Morales88"
$ws.Range("F14").Value = "D@Morales88"
$ws.Range("G14").Value = "Yes"
$ws.Range("H14").Formula = "=""false"""

# Convert the "=""false""" helper formulas (rows 9-14) into plain static
# shared-string values, same as was done for H6/H7/H8 above.
$ws.Range("H9:H14").Copy()
$ws.Range("H9:H14").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) Misc view-state tweak captured by the diff: the active selection moves
#    from D9 to F19.
# ---------------------------------------------------------------------------
$ws.Range("F19").Select()
